# spreadsheet and tooltiptext improvement
#
# Applies:
#  - fieldnames: clear the "D:/Biocentury/Gapminder" row and the extra
#    "Google Datasearch/Kaggle/OurWorldInData" tooltip rows, keep only the
#    "Expasy" hyperlink, rename the "C:" cell to "GOOGLE", move selection.
#  - category: clear the duplicate "app1/inet1/inet2" row and the extra
#    "inet2" tooltip cells, move selection.
#  - URL: clear the "D:\/bio-century/gapminder..." row and the extra
#    tooltip cells, rename "C:\" to "www.google.de" and hyperlink it,
#    drop the now-orphaned bio-century hyperlink, move selection.
#  - color: clear the duplicate "246,194,62/...": row, drop the extra
#    "76,135,200" tooltip cells, move selection, sheet no longer active.
#  - add a new "test" sheet at the end with a header row copied from
#    "fieldnames" and a single "sdfaf" value, make it the active sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "fieldnames"
# ---------------------------------------------------------------------
$fieldnames = $wb.Worksheets.Item("fieldnames")

$fieldnames.Range("A2").Value = "GOOGLE"

$fieldnames.Range("A3:C3").ClearContents()
$fieldnames.Range("C4").ClearContents()
$fieldnames.Range("C5").ClearContents()
$fieldnames.Range("C6").ClearContents()

# keep only the "Expasy" hyperlink on C2
$fieldnames.Hyperlinks.Delete()
$fieldnames.Hyperlinks.Add($fieldnames.Range("C2"), "https://www.expasy.org/", "", "", "Expasy")

$fieldnames.Range("A3").Select()

# ---------------------------------------------------------------------
# Sheet "category"
# ---------------------------------------------------------------------
$category = $wb.Worksheets.Item("category")

$category.Range("A3:C3").ClearContents()
$category.Range("C4").ClearContents()
$category.Range("C5").ClearContents()
$category.Range("C6").ClearContents()

$category.Range("F11").Select()

# ---------------------------------------------------------------------
# Sheet "URL"
# ---------------------------------------------------------------------
$url = $wb.Worksheets.Item("URL")

$url.Range("A2").Value = "www.google.de"

$url.Range("A3:C3").ClearContents()
$url.Range("C4").ClearContents()
$url.Range("C5").ClearContents()
$url.Range("C6").ClearContents()

# rebuild hyperlinks: add google, keep github + router ip, drop bio-century
$url.Hyperlinks.Delete()
$url.Hyperlinks.Add($url.Range("A2"), "https://www.google.de/", "", "", "www.google.de")
$url.Hyperlinks.Add($url.Range("B2"), "http://www.github.com/", "", "", "www.github.com")
$url.Hyperlinks.Add($url.Range("D2"), "https://192.168.178.1/", "", "", "https://192.168.178.1/")

$url.Range("A3").Select()

# ---------------------------------------------------------------------
# Sheet "color"
# ---------------------------------------------------------------------
$color = $wb.Worksheets.Item("color")

$color.Range("A3:C3").ClearContents()
$color.Range("C4").ClearContents()
$color.Range("C5").ClearContents()
$color.Range("C6").ClearContents()

$color.Range("D2").Select()

# ---------------------------------------------------------------------
# New sheet "test"
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$test = $wb.Worksheets.Add($null, $lastSheet)
$test.Name = "test"

$fieldnames.Range("A1:D1").Copy($test.Range("A1:D1"))
$fieldnames.Range("A3:D3").Copy()
$test.Range("A2:D2").PasteSpecial(-4122)
$test.Range("A2").Value = "sdfaf"

$test.Activate()
$test.Range("B6").Select()
